$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3624
$ws.Range("I19").Value = 7751.7144
$ws.Range("J19").Value = 734.6
$ws.Range("K19").Value = 7751.7144
$ws.Range("L19").Value = 734.6
$ws.Range("M19").Value = -7576.7144
$ws.Range("N19").Value = -1084.6

$ws.Range("H80").Value = 27630.611
$ws.Range("I80").Value = 56898
$ws.Range("J80").Value = 12996.917
$ws.Range("K80").Value = 170694
$ws.Range("L80").Value = 38990.751
$ws.Range("M80").Value = -169696
$ws.Range("N80").Value = -40986.751

$ws.Range("H83").Value = 27630.611
$ws.Range("I83").Value = 56898
$ws.Range("J83").Value = 12996.917
$ws.Range("K83").Value = 512082
$ws.Range("L83").Value = 116972.253
$ws.Range("M83").Value = -507090
$ws.Range("N83").Value = -126956.253

$ws.Range("H98").Value = 1749.9
$ws.Range("I98").Value = 1855.4445
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 1855.4445
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = -357.4445000000001
$ws.Range("N98").Value = -3796

$ws.Range("H107").Value = 210.1923
$ws.Range("I107").Value = 191.59091
$ws.Range("K107").Value = 191.59091
$ws.Range("M107").Value = 1728.40909

$ws.Range("H112").Value = 1846.4706
$ws.Range("I112").Value = 450
$ws.Range("J112").Value = 2032.6666
$ws.Range("K112").Value = 1350
$ws.Range("L112").Value = 6097.9998
$ws.Range("M112").Value = -242
$ws.Range("N112").Value = -8313.9998

$ws.Range("H113").Value = 3421.9375
$ws.Range("I113").Value = 2753.2856
$ws.Range("K113").Value = 2753.2856
$ws.Range("M113").Value = 500.7143999999998

$ws.Range("H116").Value = 91977
$ws.Range("I116").Value = 135273
$ws.Range("J116").Value = 5385
$ws.Range("K116").Value = 135273
$ws.Range("L116").Value = 5385
$ws.Range("M116").Value = -131831
$ws.Range("N116").Value = -12269

$ws.Range("H122").Value = 1749.9
$ws.Range("I122").Value = 1855.4445
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 5566.333500000001
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -3116.333500000001
$ws.Range("N122").Value = -7300

$ws.Range("H125").Value = 626.1667
$ws.Range("I125").Value = 687.1579
$ws.Range("J125").Value = 394.4
$ws.Range("K125").Value = 6184.4211
$ws.Range("L125").Value = 3549.6
$ws.Range("M125").Value = -3724.4211
$ws.Range("N125").Value = -8469.6

$ws.Range("H137").Value = 288488.47
$ws.Range("I137").Value = 2357.5144
$ws.Range("K137").Value = 7072.5432
$ws.Range("M137").Value = -4522.5432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3021.03
$ws.Range("I32").Value = 3021.03
$ws.Range("K32").Value = 3021.03
$ws.Range("M32").Value = -2734.03

$ws.Range("H97").Value = 719.7826
$ws.Range("I97").Value = 509.72223
$ws.Range("K97").Value = 509.72223
$ws.Range("M97").Value = -13.72223000000002

$ws.Range("H110").Value = 2360.2144
$ws.Range("I110").Value = 2218.6924
$ws.Range("K110").Value = 2218.6924
$ws.Range("M110").Value = -173.6923999999999

$ws.Range("H132").Value = 3533.238
$ws.Range("I132").Value = 3929.3044
$ws.Range("J132").Value = 3053.7896
$ws.Range("K132").Value = 11787.9132
$ws.Range("L132").Value = 9161.3688
$ws.Range("M132").Value = -9257.913199999999
$ws.Range("N132").Value = -14221.3688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2309.6538
$ws.Range("I20").Value = 2202.4443
$ws.Range("J20").Value = 2550.875
$ws.Range("K20").Value = 2202.4443
$ws.Range("L20").Value = 2550.875
$ws.Range("M20").Value = -1955.4443
$ws.Range("N20").Value = -3044.875

$ws.Range("H134").Value = 20435.701
$ws.Range("I134").Value = 26635.875
$ws.Range("J134").Value = 5847.0586
$ws.Range("K134").Value = 79907.625
$ws.Range("L134").Value = 17541.1758
$ws.Range("M134").Value = -77372.625
$ws.Range("N134").Value = -22611.1758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11629697
$ws.Range("I134").Value = 20001234
$ws.Range("J134").Value = 2561.389
$ws.Range("K134").Value = 60003702
$ws.Range("L134").Value = 7684.167
$ws.Range("M134").Value = -60001167
$ws.Range("N134").Value = -12754.167

$ws.Range("H140").Value = 15555
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 35.884617
$ws.Range("J12").Value = 36.666668
$ws.Range("L12").Value = 110.000004
$ws.Range("N12").Value = -456.000004

$ws.Range("H129").Value = 1577.1613
$ws.Range("I129").Value = 739
$ws.Range("J129").Value = 1976.2858
$ws.Range("K129").Value = 2217
$ws.Range("L129").Value = 5928.857400000001
$ws.Range("M129").Value = 2783
$ws.Range("N129").Value = -15928.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4168.9062
$ws.Range("I80").Value = 4551.522
$ws.Range("J80").Value = 3191.111
$ws.Range("K80").Value = 4551.522
$ws.Range("L80").Value = 3191.111
$ws.Range("M80").Value = -3553.522
$ws.Range("N80").Value = -5187.111

$ws.Range("H83").Value = 4168.9062
$ws.Range("I83").Value = 4551.522
$ws.Range("J83").Value = 3191.111
$ws.Range("K83").Value = 22757.61
$ws.Range("L83").Value = 15955.555
$ws.Range("M83").Value = -17765.61
$ws.Range("N83").Value = -25939.555

$ws.Range("H102").Value = 3805.75
$ws.Range("I102").Value = 4227.579
$ws.Range("J102").Value = 2202.8
$ws.Range("K102").Value = 4227.579
$ws.Range("L102").Value = 2202.8
$ws.Range("M102").Value = -2605.579
$ws.Range("N102").Value = -5446.8

$ws.Range("H126").Value = 3998.8
$ws.Range("I126").Value = 4724.75
$ws.Range("J126").Value = 3514.8333
$ws.Range("K126").Value = 14174.25
$ws.Range("L126").Value = 10544.4999
$ws.Range("M126").Value = -11704.25
$ws.Range("N126").Value = -15484.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2207.2083
$ws.Range("I7").Value = 2256.5833
$ws.Range("J7").Value = 2157.8333
$ws.Range("K7").Value = 2256.5833
$ws.Range("L7").Value = 2157.8333
$ws.Range("M7").Value = -2144.5833
$ws.Range("N7").Value = -2381.8333

$ws.Range("H122").Value = 2968.5789
$ws.Range("I122").Value = 3140.4
$ws.Range("J122").Value = 2777.6667
$ws.Range("K122").Value = 9421.200000000001
$ws.Range("L122").Value = 8333.000100000001
$ws.Range("M122").Value = -6971.200000000001
$ws.Range("N122").Value = -13233.0001

$ws.Range("H126").Value = 2207.2083
$ws.Range("I126").Value = 2256.5833
$ws.Range("J126").Value = 2157.8333
$ws.Range("K126").Value = 6769.749899999999
$ws.Range("L126").Value = 6473.499899999999
$ws.Range("M126").Value = -4299.749899999999
$ws.Range("N126").Value = -11413.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 64888.5
$ws.Range("I122").Value = 73815.42999999999
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 221446.29
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -218996.29
$ws.Range("N122").Value = -12100

$ws.Range("H126").Value = 1144.4375
$ws.Range("I126").Value = 938.5217
$ws.Range("J126").Value = 1670.6666
$ws.Range("K126").Value = 2815.5651
$ws.Range("L126").Value = 5011.9998
$ws.Range("M126").Value = -345.5650999999998
$ws.Range("N126").Value = -9951.9998
